# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were re-derived from the source data and
# need to be overwritten with the newly computed strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new K value (column G), per the regenerated save_data.
$kValues = @{
    2  = 4
    4  = 1
    5  = 2
    6  = 2
    7  = 2
    8  = 3
    9  = 2
    10 = 2
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 2
    21 = 1
    22 = 2
    23 = 0
    24 = 1
    25 = 1
    26 = 1
    27 = 2
    28 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
